# Update the latest test run timestamps (column B) and result (column A)
# for the newly bootstrapped "BWP" test sheets, per the Katalon run logs.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("CCDeferredPS_27")
$ws.Range("B2").Value = "Sat Nov 08 13:48:21 IST 2025"

$ws = $wb.Worksheets.Item("CCDeferredPC_27")
$ws.Range("B2").Value = "Sat Nov 08 13:46:08 IST 2025"

$ws = $wb.Worksheets.Item("CCDeferredCorp_27")
$ws.Range("B2").Value = "Sat Nov 08 13:44:53 IST 2025"

$ws = $wb.Worksheets.Item("CMCAutopayPC_27")
$ws.Range("B2").Value = "Sat Nov 08 13:51:40 IST 2025"

$ws = $wb.Worksheets.Item("CMCAutopayCorp_27")
$ws.Range("B2").Value = "Sat Nov 08 13:50:42 IST 2025"

$ws = $wb.Worksheets.Item("CMCAutopayPS_27")
$ws.Range("B2").Value = "Sat Nov 08 13:52:37 IST 2025"

$ws = $wb.Worksheets.Item("PayNowDCFCorp_27")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Sat Nov 08 13:56:05 IST 2025"
